$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 2020
$ws.Range("B18").Value = 0.3

$ws.Range("A19").Select()
